$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.782.13'
$ws.Range("E2").Value = '  -1.46%  '

# Row 3
$ws.Range("D3").Value = '2.185.71'
$ws.Range("E3").Value = '  -2.70%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.13'
$ws.Range("E5").Value = '  -3.99%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '88.12'
$ws.Range("E6").Value = '  -5.55%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.561'
$ws.Range("E7").Value = '  -1.64%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.477'
$ws.Range("E9").Value = '  -8.70%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.27'
$ws.Range("E10").Value = '  -6.69%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0765'
$ws.Range("E11").Value = '  -5.70%  '

# Row 12
$ws.Range("E12").Value = '  -1.83%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.74'
$ws.Range("E13").Value = '  -5.71%  '

# Row 14
$ws.Range("D14").Value = '2.519.41'
$ws.Range("E14").Value = '  -2.61%  '

# Row 15
$ws.Range("D15").Value = '2.266.79'
$ws.Range("E15").Value = '  +1.49%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '12.92'
$ws.Range("E16").Value = '  -4.60%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.760'
$ws.Range("E17").Value = '  -9.09%  '

# Row 18
$ws.Range("D18").Value = '43.400.91'
$ws.Range("E18").Value = '  -1.64%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0877'
$ws.Range("E19").Value = '  -8.79%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.81'
$ws.Range("E20").Value = '  -8.50%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.78'
$ws.Range("E21").Value = '  -12.42%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '62.66'
$ws.Range("E22").Value = '  -4.56%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '227.19'
$ws.Range("E23").Value = '  -4.42%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.75'
$ws.Range("E24").Value = '  -6.52%  '

# Row 25
$ws.Range("E25").Value = '  -0.04%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.80'
$ws.Range("E26").Value = '  -9.04%  '

# Row 27
$ws.Range("E27").Value = '  -0.63%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.14'
$ws.Range("E28").Value = '  -6.80%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.23'
$ws.Range("E29").Value = '  -8.53%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.94'
$ws.Range("E30").Value = '  -5.59%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '145.85'
$ws.Range("E31").Value = '  -4.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.26'
$ws.Range("E32").Value = '  -11.33%  '

# Row 33
$ws.Range("E33").Value = '  -7.18%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0721'
$ws.Range("E34").Value = '  -9.71%  '

# Row 35
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.89'
$ws.Range("E35").Value = '  -6.82%  '

# Row 36
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.115'
$ws.Range("E36").Value = '  -3.98%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.101'
$ws.Range("E37").Value = '  -7.93%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.63'
$ws.Range("E38").Value = '  -9.26%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0278'
$ws.Range("E39").Value = '  -7.50%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.47'
$ws.Range("E40").Value = '  -9.30%  '

# Row 41
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.04'
$ws.Range("E41").Value = '  -12.31%  '

# Row 42
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.10'
$ws.Range("E42").Value = '  -10.97%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.01'
$ws.Range("E43").Value = '  -0.22%  '

# Row 44
$ws.Range("D44").Value = '1.746.82'
$ws.Range("E44").Value = '  +0.79%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.64'
$ws.Range("E45").Value = '  +1.24%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.83'
$ws.Range("E46").Value = '  -1.57%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '72.93'
$ws.Range("E47").Value = '  -9.55%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.170'
$ws.Range("E48").Value = '  -11.47%  '

# Row 49
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '90.85'
$ws.Range("E49").Value = '  -8.77%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.402.92'
$ws.Range("E50").Value = '  -2.56%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.37'
